$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# The "N" / "safari" row (row 5) of the HomePageTest block is no longer
# needed, so remove it entirely (shifts everything below up by one row).
$ws.Rows.Item(5).Delete()

# The remaining "firefox" run in the HomePageTest block should log as
# "chrome" instead (each data row now logs its own result separately).
$ws.Cells.Item(4, 2).Value = "chrome"

# Same fix for the BuyBodyLotionTest block: the last data row (now row 9
# after the deletion above) also switches from "firefox" to "chrome".
$ws.Cells.Item(9, 2).Value = "chrome"

# Hyperlinks are anchored by cell reference and do not follow the row
# shift automatically, so re-create them pointing at the correct cells.
$ws.Range("G8").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G8"), "mailto:swapbamnote@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G9"), "mailto:swapbamnote@gmail.com") | Out-Null
$ws.Range("G8").Style = "Hyperlink"
$ws.Range("G9").Style = "Hyperlink"

# Reflect the final selection left on the sheet.
$ws.Range("A9:K9").Select() | Out-Null
